$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("01-RegistrarAtencion")

# Clear the "X" marker cells in column A, rows 2-6 on sheet 1
$ws1.Range("A2:A6").Clear()

# Update the placa codes in column D (rows 2-10) on sheet 1
$ws1.Range("D2").Value = "ABA259"
$ws1.Range("D3").Value = "ABA260"
$ws1.Range("D4").Value = "ABA261"
$ws1.Range("D5").Value = "ABA262"
$ws1.Range("D6").Value = "ABA263"
$ws1.Range("D7").Value = "ABA264"
$ws1.Range("D8").Value = "ABA265"
$ws1.Range("D9").Value = "ABA266"
$ws1.Range("D10").Value = "ABA267"

# Make sheet 1 the active sheet/tab and set its selection to D2:D10 (active cell D2)
$ws1.Activate()
$ws1.Range("D2:D10").Select()
